# Update risk_to_assets (B), resilience (C) and risk (D) category values
# to match the refreshed model/data snapshot (as of March 17th).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tiers")

$ws.Range("B2").Value = "Mid"
$ws.Range("D2").Value = "Mid"
$ws.Range("B4").Value = "Mid"
$ws.Range("D4").Value = "Mid"
$ws.Range("B5").Value = "High"
$ws.Range("B7").Value = "High"
$ws.Range("D7").Value = "High"
$ws.Range("B8").Value = "Mid"
$ws.Range("D8").Value = "Mid"
$ws.Range("B13").Value = "High"
$ws.Range("D13").Value = "High"
$ws.Range("D14").Value = "Low"
$ws.Range("B16").Value = "Mid"
$ws.Range("C16").Value = "High"
$ws.Range("D16").Value = "Mid"
$ws.Range("B18").Value = "Mid"
$ws.Range("D18").Value = "Mid"
$ws.Range("B19").Value = "High"
$ws.Range("D19").Value = "High"
$ws.Range("B20").Value = "Mid"
$ws.Range("D20").Value = "Mid"
$ws.Range("B21").Value = "High"
$ws.Range("D21").Value = "High"
$ws.Range("B22").Value = "High"
$ws.Range("D22").Value = "Mid"
$ws.Range("B23").Value = "High"
$ws.Range("D23").Value = "High"
$ws.Range("B24").Value = "Mid"
$ws.Range("D24").Value = "Mid"
$ws.Range("B25").Value = "Mid"
$ws.Range("D25").Value = "High"
$ws.Range("B26").Value = "High"
$ws.Range("B28").Value = "Low"
$ws.Range("D28").Value = "Low"
$ws.Range("B29").Value = "Low"
$ws.Range("B30").Value = "Low"
$ws.Range("B31").Value = "Low"
$ws.Range("B32").Value = "Low"
$ws.Range("B37").Value = "High"
$ws.Range("D37").Value = "High"
$ws.Range("B38").Value = "High"
$ws.Range("D39").Value = "Mid"
$ws.Range("D40").Value = "High"
$ws.Range("B41").Value = "Low"
$ws.Range("D41").Value = "Mid"
$ws.Range("D42").Value = "Low"
$ws.Range("B44").Value = "Mid"
$ws.Range("D44").Value = "Mid"
$ws.Range("B45").Value = "Mid"
$ws.Range("D45").Value = "Mid"
$ws.Range("B46").Value = "Low"
$ws.Range("B49").Value = "Mid"
$ws.Range("D49").Value = "Low"
$ws.Range("B50").Value = "Low"
$ws.Range("D50").Value = "Mid"
$ws.Range("B51").Value = "Low"
$ws.Range("D51").Value = "Low"
$ws.Range("B53").Value = "High"
$ws.Range("D54").Value = "Low"
$ws.Range("B56").Value = "High"
$ws.Range("D56").Value = "High"
$ws.Range("B58").Value = "Mid"
$ws.Range("D58").Value = "Mid"
$ws.Range("B59").Value = "Mid"
$ws.Range("B60").Value = "Mid"
$ws.Range("D60").Value = "Low"
$ws.Range("B61").Value = "Mid"
$ws.Range("D61").Value = "Mid"
$ws.Range("B62").Value = "Mid"
$ws.Range("D62").Value = "Mid"
